$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values remain text (matches original inlineStr formatting)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.060.48'
$ws.Range("E2").Value = '  -1.72%  '
$ws.Range("D3").Value = '2.103.50'
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("E4").Value = '  -0.61%  '
$ws.Range("D5").Value = '347.55'
$ws.Range("E5").Value = '  +3.12%  '
$ws.Range("E6").Value = '  -0.58%  '
$ws.Range("D7").Value = '0.5171'
$ws.Range("E7").Value = '  -1.43%  '
$ws.Range("E8").Value = '  -2.75%  '
$ws.Range("D9").Value = '52.33'
$ws.Range("E9").Value = '  -3.71%  '
$ws.Range("D10").Value = '0.08960'
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").Value = '25.51'
$ws.Range("E12").Value = '  +3.67%  '
$ws.Range("D13").Value = '2.113.23'
$ws.Range("E13").Value = '  -0.45%  '
$ws.Range("D14").Value = '8.229'
$ws.Range("E14").Value = '  +1.34%  '
$ws.Range("D15").Value = '6.726'
$ws.Range("E15").Value = '  -2.10%  '
$ws.Range("D16").Value = '99.38'
$ws.Range("E16").Value = '  +2.22%  '
$ws.Range("E17").Value = '  -2.08%  '
$ws.Range("E18").Value = '  -0.52%  '
$ws.Range("D19").Value = '20.81'
$ws.Range("E19").Value = '  +7.07%  '
$ws.Range("D20").Value = '0.06679'
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("E21").Value = '  -0.54%  '
$ws.Range("D22").Value = '6.238'
$ws.Range("E22").Value = '  -1.20%  '
$ws.Range("D23").Value = '30.158.18'
$ws.Range("E23").Value = '  -1.63%  '
$ws.Range("D24").Value = '12.72'
$ws.Range("E24").Value = '  -1.09%  '
$ws.Range("D25").Value = '2.346'
$ws.Range("D26").Value = '2.361.29'
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("D27").Value = '21.96'
$ws.Range("E27").Value = '  -2.05%  '
$ws.Range("D28").Value = '2.532'
$ws.Range("E28").Value = '  -0.69%  '
$ws.Range("D29").Value = '162.25'
$ws.Range("E29").Value = '  -1.11%  '
$ws.Range("D30").Value = '133.67'
$ws.Range("E30").Value = '  -0.37%  '
$ws.Range("D31").Value = '1.176'
$ws.Range("E31").Value = '  -3.19%  '
$ws.Range("D33").Value = '1.636'
$ws.Range("E33").Value = '  -1.17%  '
$ws.Range("D34").Value = '6.236'
$ws.Range("E34").Value = '  -2.19%  '
$ws.Range("D35").Value = '3.958'
$ws.Range("E35").Value = '  +0.17%  '
$ws.Range("D36").Value = '5.917'
$ws.Range("E36").Value = '  +0.30%  '
$ws.Range("D37").Value = '10.23'
$ws.Range("E37").Value = '  -3.53%  '
$ws.Range("D38").Value = '0.02572'
$ws.Range("E38").Value = '  -1.97%  '
$ws.Range("D39").Value = '0.06803'
$ws.Range("E39").Value = '  -0.58%  '
$ws.Range("E40").Value = '  -1.47%  '
$ws.Range("D41").Value = '12.56'
$ws.Range("E41").Value = '  -0.60%  '
$ws.Range("D42").Value = '0.6812'
$ws.Range("E42").Value = '  -1.16%  '
$ws.Range("E43").Value = '  +2.22%  '
$ws.Range("E44").Value = '  -4.05%  '
$ws.Range("D45").Value = '0.6375'
$ws.Range("E45").Value = '  -1.40%  '
$ws.Range("D46").Value = '2.291'
$ws.Range("E46").Value = '  -1.44%  '
$ws.Range("D47").Value = '0.00000000362'
$ws.Range("E47").Value = '  -1.23%  '
$ws.Range("E48").Value = '  -1.20%  '
$ws.Range("D49").Value = '1.221'
$ws.Range("E49").Value = '  -2.76%  '
$ws.Range("D50").Value = '82.30'
$ws.Range("E50").Value = '  -1.44%  '
$ws.Range("D51").Value = '0.07232'
$ws.Range("E51").Value = '  +0.41%  '
